# Underwater_Robot_NFC balance workbook update
# - refresh data values in column B
# - add four new rows (No.15 .. No.18)
# - give the chart a title ("Team Balance")
# - extend the chart source ranges to the new data
# - drop the explicit bar overlap override
# - reposition/resize the chart
# - update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing values in column B ----
$ws.Range("B1").Value  = 500
$ws.Range("B2").Value  = 500
$ws.Range("B3").Value  = 500
$ws.Range("B4").Value  = 500
$ws.Range("B5").Value  = 500
$ws.Range("B10").Value = 500
$ws.Range("B12").Value = 490
$ws.Range("B13").Value = 500
$ws.Range("B14").Value = 495

# ---- Add four new rows of data ----
$ws.Range("A15").Value = "No.15"
$ws.Range("B15").Value = 500
$ws.Range("A16").Value = "No.16"
$ws.Range("B16").Value = 500
$ws.Range("A17").Value = "No.17"
$ws.Range("B17").Value = 500
$ws.Range("A18").Value = "No.18"
$ws.Range("B18").Value = 500

# ---- Chart updates ----
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# Give the chart a title
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Team Balance"

# Extend the series ranges to cover the new rows
$ser = $chart.SeriesCollection(1)
$ser.Formula = '=SERIES(,Sheet1!$A$1:$A$18,Sheet1!$B$1:$B$18,1)'

# Remove the explicit bar-overlap override (back to the implicit default)
$chart.ChartGroups(1).Overlap = 0

# Reposition / resize the chart object on the sheet
$co.Left = 261.95
$co.Top = 16.2
$co.Width = 790.525
$co.Height = 272.4

# ---- Update the saved selection ----
$ws.Range("D15").Select() | Out-Null
